# orm working for postgreSQL and MySQL
#
# Applies the recorded edits:
#  - Functions sheet (sheet1): is_locked column (P) switches from text "no"/"yes"
#    to numeric 0/1; a new "message" note is added at M9; Functions becomes the
#    active/selected sheet with selection on M14.
#  - function_parameters sheet (sheet3): adds numeric ids (1,2,3,4) to the
#    previously-empty D column for the "text"-typed parameter rows; selection
#    moves to D7.
#  - Groups sheet (sheet9) is no longer the tab-selected sheet (Functions is).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Functions sheet
# ---------------------------------------------------------------------------
$wsFunctions = $wb.Worksheets.Item("Functions")

# is_locked (P) column: "no" -> 0, "yes" -> 1 (text -> numeric boolean flag)
$wsFunctions.Range("P2").Value = 0
$wsFunctions.Range("P3").Value = 1
$wsFunctions.Range("P4").Value = 0
$wsFunctions.Range("P5").Value = 1
$wsFunctions.Range("P6").Value = 0
$wsFunctions.Range("P7").Value = 1
$wsFunctions.Range("P8").Value = 0
$wsFunctions.Range("P9").Value = 1

# New version_comments note for row 9
$wsFunctions.Range("M9").Value = "message"

# ---------------------------------------------------------------------------
# function_parameters sheet
# ---------------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("function_parameters")

$wsParams.Range("D2").Value = 1
$wsParams.Range("D3").Value = 2
$wsParams.Range("D5").Value = 3
$wsParams.Range("D6").Value = 4

# ---------------------------------------------------------------------------
# View state: Functions becomes the active/selected sheet (was Groups)
# ---------------------------------------------------------------------------
$wsFunctions.Activate() | Out-Null
$wsFunctions.Range("M14").Select() | Out-Null

$wsParams.Range("D7").Select() | Out-Null

$wsGroups = $wb.Worksheets.Item("Groups")
$wsGroups.Range("B1").Select() | Out-Null

$wsFunctions.Activate() | Out-Null
